$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Stinkers" header text from G2, leaving the cell blank
# (keeps its style, but no string value/type)
$ws.Range("G2").ClearContents()

# Move the active selection from D2 to G2
$ws.Range("G2").Select()
